# Financials update: add a new quarter (two new columns) ahead of the
# existing quarterly data on the FMBM sheet, shifting the old data right,
# and touch up a handful of previously-reported figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two new columns before column D (the first data column). ---
# Everything that used to live in D:K now lives in F:M.
$ws.Range("D1:E1").EntireColumn.Insert()

# Pick up the number formatting (date format for the "Period Ending" rows,
# #,##0 for the data rows) from the column immediately to the right (the
# old column D, now column F) so the two new columns look like the rest of
# the table instead of using the default "General" style. Restrict the
# copy to the sheet's used rows so we don't balloon the used range out to
# the full 1,048,576 rows of the worksheet.
$lastRow = $ws.UsedRange.Rows($ws.UsedRange.Rows.Count).Row
$ws.Range("F5:F" + $lastRow).Copy()
$ws.Range("D5:E" + $lastRow).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Columns("D:E").ColumnWidth = $ws.Columns("F").ColumnWidth

# --- 2. Populate the two new columns with the new quarter's figures. ---
# Row number -> @(column D value, column E value)
$newQuarterData = @{
    7   = @(43465, 43373)
    8   = @(9700, 9300)
    9   = @("NA", "NA")
    10  = @("NA", "NA")
    12  = @("NA", "NA")
    13  = @(0, 0)
    14  = @(0, 0)
    15  = @(0, 0)
    17  = @(1900, 1800)
    18  = @(7800, 7500)
    20  = @(-4600, -4800)
    21  = @(3600, 3100)
    22  = @(0, 0)
    23  = @(3300, 2800)
    24  = @(300, 300)
    25  = @(0, 0)
    26  = @(2900, 2500)
    27  = @(2800, 2400)
    28  = @(0, 0)
    29  = @(0, "NA")
    30  = @(0, 0)
    31  = @(0, 0)
    32  = @(4600, 4800)
    33  = @(2800, 2400)
    34  = @(0, 0)
    35  = @(2800, 2400)
    38  = @(43465, 43373)
    41  = @(10900, 11100)
    42  = @(13400, 20200)
    43  = @(0, 0)
    44  = @(0, 0)
    45  = @(0, 0)
    46  = @(0, 0)
    47  = @(0, 0)
    48  = @(17800, 17600)
    49  = @(2900, 3000)
    50  = @(0, 0)
    51  = @(0, 0)
    52  = @(0, 0)
    53  = @(0, 0)
    54  = @(780300, 775600)
    57  = @(0, 0)
    58  = @(0, 0)
    59  = @(16700, 17600)
    60  = @(0, 0)
    61  = @(0, 0)
    62  = @(0, 0)
    63  = @(0, 0)
    64  = @(0, 0)
    65  = @(0, 0)
    66  = @(688900, 683500)
    68  = @(0, 0)
    69  = @(0, 0)
    70  = @(5700, 7500)
    71  = @(0, 0)
    72  = @(65600, 63500)
    73  = @(0, 0)
    74  = @(0, 0)
    75  = @(0, 0)
    76  = @(85700, 84600)
    77  = @(0, 0)
    80  = @(43465, 43373)
    81  = @(2800, 2400)
    83  = @(300, 300)
    84  = @(0, 0)
    85  = @(0, 0)
    86  = @(0, 0)
    87  = @(0, 0)
    88  = @(0, 0)
    89  = @(3600, 3700)
    91  = @(-500, -400)
    92  = @(0, 0)
    93  = @(0, 0)
    94  = @(-12400, -1400)
    96  = @(-900, -900)
    97  = @(0, 0)
    98  = @(0, 0)
    99  = @(0, 0)
    100 = @(2300, 1700)
    101 = @(0, 0)
    102 = @(-6500, 4100)
}

foreach ($row in $newQuarterData.Keys) {
    $vals = $newQuarterData[$row]
    $ws.Range("D$row").Value2 = $vals[0]
    $ws.Range("E$row").Value2 = $vals[1]
}

# --- 3. A few previously-reported quarters were also revised. ---
# "Other Assets" (row 52): cells that used to read "NA" now read 0.
foreach ($col in @("F", "G", "H", "I", "J")) {
    $ws.Range($col + "52").Value2 = 0
}

# "Total Cash Flow From Operating Activities" (row 89): restated figure.
$ws.Range("H89").Value2 = -200

# "Capital Expenditures" (row 91): restated figure.
$ws.Range("H91").Value2 = 3000
